$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 393, pushing existing rows 393:458 down to 394:459
$ws.Rows.Item(393).Insert()

# Columns that stay constant across the table (copy down from the row now below, i.e. row 394)
$ws.Range("A393").Value = $ws.Range("A394").Value2
$ws.Range("B393").Value = $ws.Range("B394").Value2
$ws.Range("C393").Value = $ws.Range("C394").Value2
$ws.Range("E393").Value = $ws.Range("E394").Value2
$ws.Range("F393").Value = $ws.Range("F394").Value2
$ws.Range("G393").Value = $ws.Range("G394").Value2
$ws.Range("H393").Value = $ws.Range("H394").Value2
$ws.Range("I393").Value = $ws.Range("I394").Value2
$ws.Range("R393").Value = $ws.Range("R394").Value2

# New weekly record values
$ws.Range("D393").Value = 45218
$ws.Range("J393").Value = 80
$ws.Range("K393").Value = 6000
$ws.Range("L393").Value = 6000
$ws.Range("M393").Value = 6000
$ws.Range("N393").Value = "$/docena de atados (3 kilos)"
$ws.Range("O393").Value = "Región Metropolitana"
$ws.Range("P393").Value = 2000
$ws.Range("Q393").Value = 3
